# Scheduled runner update: refresh market-price / profit figures across
# the per-job Sheets (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1404
$ws.Range("I32").Value = 1057.2858
$ws.Range("J32").Value = 1538.8334
$ws.Range("K32").Value = 1057.2858
$ws.Range("L32").Value = 1538.8334
$ws.Range("M32").Value = -731.2858000000001
$ws.Range("N32").Value = -2190.8334

$ws.Range("H116").Value = 5807.7036
$ws.Range("I116").Value = 6141.273
$ws.Range("J116").Value = 4340
$ws.Range("K116").Value = 6141.273
$ws.Range("L116").Value = 4340
$ws.Range("M116").Value = -2699.273
$ws.Range("N116").Value = -11224

$ws.Range("H138").Value = 2867.861
$ws.Range("I138").Value = 761.3333
$ws.Range("J138").Value = 5817
$ws.Range("K138").Value = 2283.9999
$ws.Range("L138").Value = 17451
$ws.Range("M138").Value = 2856.0001
$ws.Range("N138").Value = -27731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4510.896
$ws.Range("I32").Value = 3065.3936
$ws.Range("J32").Value = 10021.875
$ws.Range("K32").Value = 3065.3936
$ws.Range("L32").Value = 10021.875
$ws.Range("M32").Value = -2778.3936
$ws.Range("N32").Value = -10595.875

$ws.Range("H110").Value = 406.875
$ws.Range("I110").Value = 440.33334
$ws.Range("J110").Value = 306.5
$ws.Range("K110").Value = 440.33334
$ws.Range("L110").Value = 306.5
$ws.Range("M110").Value = 1604.66666
$ws.Range("N110").Value = -4396.5

$ws.Range("H122").Value = 1605737
$ws.Range("I122").Value = 2568356.8
$ws.Range("J122").Value = 1370.6666
$ws.Range("K122").Value = 7705070.399999999
$ws.Range("L122").Value = 4111.9998
$ws.Range("M122").Value = -7702620.399999999
$ws.Range("N122").Value = -9011.9998

$ws.Range("H123").Value = 37999.5
$ws.Range("J123").Value = 37999.5
$ws.Range("L123").Value = 37999.5
$ws.Range("N123").Value = -47799.5

$ws.Range("H137").Value = 49978.184
$ws.Range("J137").Value = 49978.184
$ws.Range("L137").Value = 49978.184
$ws.Range("N137").Value = -60178.184

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws.Range("H139").Value = 48000
$ws.Range("J139").Value = 48000
$ws.Range("L139").Value = 48000
$ws.Range("N139").Value = -58280

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7329.2085
$ws.Range("J4").Value = 5586.364
$ws.Range("L4").Value = 5586.364
$ws.Range("N4").Value = -5810.364

$ws.Range("H31").Value = 8774625
$ws.Range("I31").Value = 1194.3658
$ws.Range("J31").Value = 31256542
$ws.Range("K31").Value = 1194.3658
$ws.Range("L31").Value = 31256542
$ws.Range("M31").Value = -899.3658
$ws.Range("N31").Value = -31257132

$ws.Range("H34").Value = 8774625
$ws.Range("I34").Value = 1194.3658
$ws.Range("J34").Value = 31256542
$ws.Range("K34").Value = 1194.3658
$ws.Range("L34").Value = 31256542
$ws.Range("M34").Value = -992.3658
$ws.Range("N34").Value = -31256946

$ws.Range("H58").Value = 4066293.2
$ws.Range("I58").Value = 6667560
$ws.Range("J58").Value = 1813.75
$ws.Range("K58").Value = 6667560
$ws.Range("L58").Value = 1813.75
$ws.Range("M58").Value = -6667357
$ws.Range("N58").Value = -2219.75

$ws.Range("H99").Value = 6583544
$ws.Range("I99").Value = 3372.4707
$ws.Range("K99").Value = 3372.4707
$ws.Range("M99").Value = -1874.4707

$ws.Range("H126").Value = 6583544
$ws.Range("I126").Value = 3372.4707
$ws.Range("K126").Value = 10117.4121
$ws.Range("M126").Value = -7647.4121

$ws.Range("H136").Value = 4066293.2
$ws.Range("I136").Value = 6667560
$ws.Range("J136").Value = 1813.75
$ws.Range("K136").Value = 20002680
$ws.Range("L136").Value = 5441.25
$ws.Range("M136").Value = -20000130
$ws.Range("N136").Value = -10541.25

$ws.Range("H138").Value = 40780
$ws.Range("J138").Value = 40780
$ws.Range("L138").Value = 40780
$ws.Range("N138").Value = -51060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 108110.71
$ws.Range("I4").Value = 150320
$ws.Range("J4").Value = 2587.5
$ws.Range("K4").Value = 450960
$ws.Range("L4").Value = 7762.5
$ws.Range("M4").Value = -450848
$ws.Range("N4").Value = -7986.5

$ws.Range("H134").Value = 10949.964
$ws.Range("I134").Value = 12950
$ws.Range("J134").Value = 9838.833000000001
$ws.Range("K134").Value = 38850
$ws.Range("L134").Value = 29516.499
$ws.Range("M134").Value = -33780
$ws.Range("N134").Value = -39656.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 600
$ws.Range("I5").Value = 600
$ws.Range("K5").Value = 600
$ws.Range("M5").Value = -488

$ws.Range("H122").Value = 4418198.5
$ws.Range("I122").Value = 4322666
$ws.Range("J122").Value = 4548469.5
$ws.Range("K122").Value = 12967998
$ws.Range("L122").Value = 13645408.5
$ws.Range("M122").Value = -12965548
$ws.Range("N122").Value = -13650308.5

$ws.Range("H132").Value = 5210677.5
$ws.Range("I132").Value = 6946388.5
$ws.Range("J132").Value = 3544.25
$ws.Range("K132").Value = 20839165.5
$ws.Range("L132").Value = 10632.75
$ws.Range("M132").Value = -20836635.5
$ws.Range("N132").Value = -15692.75

$ws.Range("H135").Value = 40350
$ws.Range("J135").Value = 40350
$ws.Range("L135").Value = 40350
$ws.Range("N135").Value = -50490

$ws.Range("H140").Value = 42931.48
$ws.Range("J140").Value = 42931.48
$ws.Range("L140").Value = 42931.48
$ws.Range("N140").Value = -53291.48

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5050

$ws.Range("H132").Value = 13100232
$ws.Range("I132").Value = 14790037
$ws.Range("K132").Value = 44370111
$ws.Range("M132").Value = -44367581

$ws.Range("H137").Value = 39700
$ws.Range("J137").Value = 39700
$ws.Range("L137").Value = 39700
$ws.Range("N137").Value = -49900

$ws.Range("H141").Value = 44213
$ws.Range("J141").Value = 44213
$ws.Range("L141").Value = 44213
$ws.Range("N141").Value = -54573

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50076.25
$ws.Range("I2").Value = 302
$ws.Range("J2").Value = 66667.664
$ws.Range("K2").Value = 302
$ws.Range("L2").Value = 66667.664
$ws.Range("M2").Value = -190
$ws.Range("N2").Value = -66891.664

$ws.Range("H107").Value = 83334010
$ws.Range("I107").Value = 111111680
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 333335040
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -333333120
$ws.Range("N107").Value = -6840
